$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'275.18"
$ws.Range("D3").Value = "'21.16"
$ws.Range("D4").Value = "'6.250"
$ws.Range("D5").Value = "'0.06184"
$ws.Range("D7").Value = "'1.519"
$ws.Range("D8").Value = "'6.535"
$ws.Range("D9").Value = "'0.8226"
$ws.Range("D10").Value = "'0.1650"
$ws.Range("D11").Value = "'0.08265"
$ws.Range("D12").Value = "'0.03457"
$ws.Range("D13").Value = "'0.03159"
$ws.Range("D14").Value = "'0.09134"
$ws.Range("D15").Value = "'3.763"
$ws.Range("D16").Value = "'0.001642"
$ws.Range("D17").Value = "'0.04673"
$ws.Range("D18").Value = "'0.006426"
$ws.Range("E18").Value = "17TigerCashTCHBestin24h"
$ws.Range("D19").Value = "'0.006133"
$ws.Range("D21").Value = "'0.0001499"
$ws.Range("D22").Value = "'3.724"
$ws.Range("D23").Value = "'2.321"
$ws.Range("D24").Value = "'0.01389"
$ws.Range("D25").Value = "'0.3326"
$ws.Range("D28").Value = "'0.0002736"
$ws.Range("D40").Value = "'0.04741"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.007037"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.004397"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("D43").Value = "'0.1105"
$ws.Range("D44").Value = "'0.01157"
$ws.Range("D45").Value = "'0.00006057"
$ws.Range("D47").Value = "'0.7229"
$ws.Range("D49").Value = "'0.00001899"
